$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 455.30768
$ws.Range("I80").Value = 361.75
$ws.Range("K80").Value = 1085.25
$ws.Range("M80").Value = -87.25
$ws.Range("H83").Value = 455.30768
$ws.Range("I83").Value = 361.75
$ws.Range("K83").Value = 3255.75
$ws.Range("M83").Value = 1736.25
$ws.Range("H113").Value = 11366355
$ws.Range("I113").Value = 2997.8572
$ws.Range("J113").Value = 31252230
$ws.Range("K113").Value = 2997.8572
$ws.Range("L113").Value = 31252230
$ws.Range("M113").Value = 256.1428000000001
$ws.Range("N113").Value = -31258738
$ws.Range("H116").Value = 6495.2085
$ws.Range("I116").Value = 8813.214
$ws.Range("J116").Value = 3250
$ws.Range("K116").Value = 8813.214
$ws.Range("L116").Value = 3250
$ws.Range("M116").Value = -5371.214
$ws.Range("N116").Value = -10134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1367.65
$ws.Range("I32").Value = 1208.1765
$ws.Range("J32").Value = 2271.3333
$ws.Range("K32").Value = 1208.1765
$ws.Range("L32").Value = 2271.3333
$ws.Range("M32").Value = -921.1765
$ws.Range("N32").Value = -2845.3333
$ws.Range("H45").Value = 6010.8
$ws.Range("I45").Value = 6762.706
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 6762.706
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -6385.706
$ws.Range("N45").Value = -2504
$ws.Range("H61").Value = 2839.8982
$ws.Range("I61").Value = 3185.9783
$ws.Range("K61").Value = 3185.9783
$ws.Range("M61").Value = -2973.9783
$ws.Range("H74").Value = 1551.72
$ws.Range("I74").Value = 1314.6842
$ws.Range("K74").Value = 1314.6842
$ws.Range("M74").Value = -440.6841999999999
$ws.Range("H77").Value = 1551.72
$ws.Range("I77").Value = 1314.6842
$ws.Range("K77").Value = 6573.420999999999
$ws.Range("M77").Value = -2205.420999999999
$ws.Range("H110").Value = 829.3158
$ws.Range("I110").Value = 736.9375
$ws.Range("J110").Value = 1322
$ws.Range("K110").Value = 736.9375
$ws.Range("L110").Value = 1322
$ws.Range("M110").Value = 1308.0625
$ws.Range("N110").Value = -5412
$ws.Range("H136").Value = 2839.8982
$ws.Range("I136").Value = 3185.9783
$ws.Range("K136").Value = 9557.9349
$ws.Range("M136").Value = -7007.9349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 198.94444
$ws.Range("I64").Value = 143.1
$ws.Range("K64").Value = 143.1
$ws.Range("M64").Value = 81.90000000000001
$ws.Range("H67").Value = 198.94444
$ws.Range("I67").Value = 143.1
$ws.Range("K67").Value = 143.1
$ws.Range("M67").Value = 636.9
$ws.Range("H107").Value = 90910450
$ws.Range("I107").Value = 333334600
$ws.Range("J107").Value = 1392.125
$ws.Range("K107").Value = 333334600
$ws.Range("L107").Value = 1392.125
$ws.Range("M107").Value = -333332680
$ws.Range("N107").Value = -5232.125
$ws.Range("H134").Value = 3694.3877
$ws.Range("I134").Value = 3932.7297
$ws.Range("J134").Value = 2959.5
$ws.Range("K134").Value = 11798.1891
$ws.Range("L134").Value = 8878.5
$ws.Range("M134").Value = -9263.1891
$ws.Range("N134").Value = -13948.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1185.0546
$ws.Range("I58").Value = 741.3333
$ws.Range("J58").Value = 2618.6155
$ws.Range("K58").Value = 741.3333
$ws.Range("L58").Value = 2618.6155
$ws.Range("M58").Value = -538.3333
$ws.Range("N58").Value = -3024.6155
$ws.Range("H132").Value = 1590.8793
$ws.Range("I132").Value = 1202.8889
$ws.Range("J132").Value = 2933.923
$ws.Range("K132").Value = 3608.6667
$ws.Range("L132").Value = 8801.769
$ws.Range("M132").Value = -1078.6667
$ws.Range("N132").Value = -13861.769
$ws.Range("H136").Value = 1185.0546
$ws.Range("I136").Value = 741.3333
$ws.Range("J136").Value = 2618.6155
$ws.Range("K136").Value = 2223.9999
$ws.Range("L136").Value = 7855.8465
$ws.Range("M136").Value = 326.0001000000002
$ws.Range("N136").Value = -12955.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2881.975
$ws.Range("I122").Value = 485
$ws.Range("K122").Value = 4365
$ws.Range("M122").Value = -1915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1859.1154
$ws.Range("I7").Value = 1676.85
$ws.Range("J7").Value = 2466.6667
$ws.Range("K7").Value = 1676.85
$ws.Range("L7").Value = 2466.6667
$ws.Range("M7").Value = -1564.85
$ws.Range("N7").Value = -2690.6667
$ws.Range("H61").Value = 1358.7368
$ws.Range("I61").Value = 1324.5294
$ws.Range("J61").Value = 1649.5
$ws.Range("K61").Value = 1324.5294
$ws.Range("L61").Value = 1649.5
$ws.Range("M61").Value = -1122.5294
$ws.Range("N61").Value = -2053.5
$ws.Range("H93").Value = 22727942
$ws.Range("I93").Value = 682
$ws.Range("J93").Value = 55556210
$ws.Range("K93").Value = 682
$ws.Range("L93").Value = 55556210
$ws.Range("M93").Value = 566
$ws.Range("N93").Value = -55558706
$ws.Range("H100").Value = 1499.8334
$ws.Range("I100").Value = 1499.8334
$ws.Range("K100").Value = 1499.8334
$ws.Range("M100").Value = -958.8334
$ws.Range("H113").Value = 1358.7368
$ws.Range("I113").Value = 1324.5294
$ws.Range("J113").Value = 1649.5
$ws.Range("K113").Value = 1324.5294
$ws.Range("L113").Value = 1649.5
$ws.Range("M113").Value = 845.4706000000001
$ws.Range("N113").Value = -5989.5
$ws.Range("H126").Value = 1859.1154
$ws.Range("I126").Value = 1676.85
$ws.Range("J126").Value = 2466.6667
$ws.Range("K126").Value = 5030.549999999999
$ws.Range("L126").Value = 7400.000100000001
$ws.Range("M126").Value = -2560.549999999999
$ws.Range("N126").Value = -12340.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 525.5
$ws.Range("I100").Value = 388.25
$ws.Range("K100").Value = 776.5
$ws.Range("M100").Value = -235.5
$ws.Range("H107").Value = 45218240
$ws.Range("I107").Value = 71429624
$ws.Range("J107").Value = 4444972.5
$ws.Range("K107").Value = 214288872
$ws.Range("L107").Value = 13334917.5
$ws.Range("M107").Value = -214286952
$ws.Range("N107").Value = -13338757.5
$ws.Range("H136").Value = 6412221
$ws.Range("I136").Value = 2047.963
$ws.Range("K136").Value = 6143.889
$ws.Range("M136").Value = -3593.889
